$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

$ws.Range("D4").Select()
